$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Drop the female_parent/male_parent columns (D:E) entirely; this shifts the
# female_plot/male_plot columns (F:G), their data AND their widths left into
# D:E.
$ws.Range("D:E").Delete()

# The shifted-in columns carried the "female_plot"/"male_plot" headers with
# them - restore the original "female_parent"/"male_parent" header text.
$ws.Range("D1").Value = "female_parent"
$ws.Range("E1").Value = "male_parent"

$ws.Range("D1").Select()
